# widened image to increase separation
# applications chart was too close to equations/logic part
#
# The roadmap diagram on slide 1 (every connector/text-box in the shape
# tree except the "Group 32" picture/diagram in the upper-left corner)
# is shifted right by 899160 EMU (70.8 pt) so it sits further away from
# the widened image. The "Date Placeholder" automatic-date field that
# appears on the slide master and all 11 slide layouts is also bumped
# from 11/25/2017 to 12/2/2017.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# EMU-per-point used by the PowerPoint object model.
$EMU_PER_PT = 12700.0

# Helper: convert a target EMU offset into the point value to assign to
# .Left/.Top so that, after the COM layer's internal float32 round-trip,
# the saved OOXML ends up with *exactly* that EMU value (plain
# emu/12700 rounds down about half the time because of float32
# precision loss, so nudge by half an EMU before dividing).
function EmuToPt($emu) {
    return ([double]$emu + 0.5) / $EMU_PER_PT
}

# Shape (by its position in $s.Shapes) -> new absolute X offset in EMU.
# Every top-level shape in the diagram moves except "Group 32" (index 4),
# whose position is unchanged in the source diff.
$moves = @(
    @{ Index = 1;  Emu = 6147457 },  # Connector: Elbow 48
    @{ Index = 2;  Emu = 5156021 },  # Connector: Elbow 80
    @{ Index = 3;  Emu = 6154264 },  # Straight Connector 64
    @{ Index = 5;  Emu = 5140224 },  # Straight Connector 37
    @{ Index = 6;  Emu = 4602667 },  # TextBox 23
    @{ Index = 7;  Emu = 4620835 },  # TextBox 24
    @{ Index = 8;  Emu = 4925386 },  # TextBox 26
    @{ Index = 9;  Emu = 4925386 },  # TextBox 27
    @{ Index = 10; Emu = 4574147 },  # TextBox 28
    @{ Index = 11; Emu = 4925386 },  # TextBox 29
    @{ Index = 12; Emu = 4648539 },  # TextBox 30
    @{ Index = 13; Emu = 5135091 },  # Connector: Elbow 45
    @{ Index = 14; Emu = 5926930 },  # TextBox 40
    @{ Index = 15; Emu = 5926930 },  # TextBox 41
    @{ Index = 16; Emu = 6852606 },  # TextBox 42
    @{ Index = 17; Emu = 6396904 },  # TextBox 43
    @{ Index = 18; Emu = 6738282 },  # TextBox 54
    @{ Index = 19; Emu = 6850704 },  # TextBox 58
    @{ Index = 20; Emu = 5376806 },  # TextBox 84
    @{ Index = 21; Emu = 5145913 },  # Connector: Elbow 88
    @{ Index = 22; Emu = 6240988 },  # TextBox 36
    @{ Index = 23; Emu = 7179385 },  # TextBox 38
    @{ Index = 24; Emu = 4054151 },  # TextBox 44
    @{ Index = 25; Emu = 5019216 },  # Connector: Elbow 5
    @{ Index = 26; Emu = 5516680 },  # Group 25
    @{ Index = 27; Emu = 3447097 },  # TextBox 53
    @{ Index = 28; Emu = 5578540 },  # TextBox 55
    @{ Index = 29; Emu = 4846450 },  # TextBox 56
    @{ Index = 30; Emu = 5318463 },  # TextBox 57
    @{ Index = 31; Emu = 5623200 },  # TextBox 59
    @{ Index = 32; Emu = 4240851 },  # TextBox 60
    @{ Index = 33; Emu = 5497374 },  # TextBox 61
    @{ Index = 34; Emu = 3640163 }   # TextBox 6
)

foreach ($m in $moves) {
    $shp = $s.Shapes.Item($m.Index)
    $shp.Left = EmuToPt $m.Emu
}

# Bump the fixed "date updated automatically" field from 11/25/2017 to
# 12/2/2017 everywhere it is cached: the slide master and every custom
# (slide) layout.
$oldDate = "11/25/2017"
$newDate = "12/2/2017"

function UpdateDatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

UpdateDatePlaceholder $p.SlideMaster.Shapes

for ($j = 1; $j -le $p.SlideMaster.CustomLayouts.Count; $j++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($j)
    UpdateDatePlaceholder $layout.Shapes
}
